$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "NewBvDIDs",
    "GB03564138",
    $null,
    "DE*850267597",
    "DE7330003759",
    "DE2070000543",
    "DE2070071908",
    "DE8170085484",
    "FI01126502",
    "LULB188712",
    "LULB185422",
    "HK0000244354",
    "US149146115L",
    "NL34140812"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    if ($null -ne $values[$i]) {
        $ws.Cells.Item($row, 1).Value = $values[$i]
    }
}

# Row 3 must hold a zero-length text string (not a truly blank cell), matching the
# source data's empty shared string. Assigning "" through .Value clears the cell
# entirely in this engine (as in real Excel), so force it through as text instead.
$ws.Cells.Item(3, 1).Formula = "=" + [char]34 + [char]34
